# The deck currently carries the "Integral" (Red Violet) theme on the
# slide master (ppt/theme/theme1.xml) and the default "Office Theme" on
# the notes master (ppt/theme/theme2.xml). The edit swaps them so the
# slide master uses the standard Office palette.
#
# VBA's RGB() packs a colour as r + g*256 + b*65536 (i.e. 0xBBGGRR) -
# reproduce that helper locally since this host doesn't expose RGB().
function RGBVal($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# msoThemeColorSchemeIndex ordering (1-based), matching a:clrScheme's
# child order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeColors = @(
    (RGBVal 0x00 0x00 0x00),   # 1  dk1      000000
    (RGBVal 0xFF 0xFF 0xFF),   # 2  lt1      FFFFFF
    (RGBVal 0x44 0x54 0x6A),   # 3  dk2      44546A
    (RGBVal 0xE7 0xE6 0xE6),   # 4  lt2      E7E6E6
    (RGBVal 0x5B 0x9B 0xD5),   # 5  accent1  5B9BD5
    (RGBVal 0xED 0x7D 0x31),   # 6  accent2  ED7D31
    (RGBVal 0xA5 0xA5 0xA5),   # 7  accent3  A5A5A5
    (RGBVal 0xFF 0xC0 0x00),   # 8  accent4  FFC000
    (RGBVal 0x44 0x72 0xC4),   # 9  accent5  4472C4
    (RGBVal 0x70 0xAD 0x47),   # 10 accent6  70AD47
    (RGBVal 0x05 0x63 0xC1),   # 11 hlink    0563C1
    (RGBVal 0x95 0x4F 0x72)    # 12 folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i - 1]
}
